# Atualización para selección de nombres
# Prefix the "res" column (F) values on the "tasas.w" sheet with "ta."
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("tasas.w")

$used = $ws.UsedRange
$lastRow = $used.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 6)   # column F = 6
    $current = $cell.Value2
    if ($current -ne $null -and $current -ne "") {
        $cell.Value2 = "ta." + $current
    }
}
